$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K2").Value = "SI"
$ws.Range("L2").Value = "NO"
$ws.Range("K3").Select()
